# Commit: "commit as of 30/3/2020"
#
# Updates the IntentSkillMappingData test workbook:
#   - "Create", "Edit" and "Delete" sheets: rename the sample values
#     K1/K2 -> Y1/Y2, NNZ/HNZ -> XNX/XNA (column order differs on the
#     "Delete" sheet) and Test -> Rest.
#   - Switches the active/selected sheet from "Create" to "Delete".

$wb = $excel.ActiveWorkbook

# --- Sheet: Create ---
$ws = $wb.Worksheets.Item("Create")
$ws.Range("D2").Value = "Y1"
$ws.Range("E2").Value = "Y2"
$ws.Range("G2").Value = "Rest"

# --- Sheet: Edit ---
$ws = $wb.Worksheets.Item("Edit")
$ws.Range("D2").Value = "Y1"
$ws.Range("E2").Value = "Y2"
$ws.Range("G2").Value = "Rest"
$ws.Range("H2").Value = "XNX"
$ws.Range("I2").Value = "XNA"

# --- Sheet: Delete ---
$ws = $wb.Worksheets.Item("Delete")
$ws.Range("D2").Value = "XNA"
$ws.Range("E2").Value = "XNX"
$ws.Range("G2").Value = "Rest"

# --- Active tab / selection changes ---
# "Create" keeps its own remembered selection (G3) but is no longer the
# active tab.
$wsCreate = $wb.Worksheets.Item("Create")
$wsCreate.Activate()
$wsCreate.Range("G3").Select()

# "Delete" becomes the active/selected tab with G2 selected.
$wsDelete = $wb.Worksheets.Item("Delete")
$wsDelete.Activate()
$wsDelete.Range("G2").Select()
